$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 123, pushing the existing
# rows 123-168 down to 125-170 (data/styles move with them).
$ws.Rows.Item(123).Insert()
$ws.Rows.Item(123).Insert()

# --- New row 123 (weekly Fruta/hortaliza entry) ---
$ws.Cells.Item(123,1).Value = 11
$ws.Cells.Item(123,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(123,3).Value = "Bíobío"
$ws.Cells.Item(123,4).Value = 44489
$ws.Cells.Item(123,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(123,5).Value = 8
$ws.Cells.Item(123,6).Value = "Fruta"
$ws.Cells.Item(123,7).Value = 100101
$ws.Cells.Item(123,8).Value = "Berries"
$ws.Cells.Item(123,9).Value = 100112025
$ws.Cells.Item(123,10).Value = "Frutilla"
$ws.Cells.Item(123,11).Value = "Sin especificar"
$ws.Cells.Item(123,12).Value = "Primera"
$ws.Cells.Item(123,13).Value = 430
$ws.Cells.Item(123,14).Value = 7500
$ws.Cells.Item(123,15).Value = 8000
$ws.Cells.Item(123,16).Value = 7767
$ws.Cells.Item(123,17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(123,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(123,19).Value = 1110
$ws.Cells.Item(123,20).Value = 7

# --- New row 124 (weekly Fruta/hortaliza entry) ---
$ws.Cells.Item(124,1).Value = 11
$ws.Cells.Item(124,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(124,3).Value = "Bíobío"
$ws.Cells.Item(124,4).Value = 44489
$ws.Cells.Item(124,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(124,5).Value = 8
$ws.Cells.Item(124,6).Value = "Fruta"
$ws.Cells.Item(124,7).Value = 100101
$ws.Cells.Item(124,8).Value = "Berries"
$ws.Cells.Item(124,9).Value = 100112025
$ws.Cells.Item(124,10).Value = "Frutilla"
$ws.Cells.Item(124,11).Value = "Sin especificar"
$ws.Cells.Item(124,12).Value = "Segunda"
$ws.Cells.Item(124,13).Value = 450
$ws.Cells.Item(124,14).Value = 6000
$ws.Cells.Item(124,15).Value = 6500
$ws.Cells.Item(124,16).Value = 6222
$ws.Cells.Item(124,17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(124,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(124,19).Value = 889
$ws.Cells.Item(124,20).Value = 7

# --- New row 169 (appended at the end, mirrors old row 168's date) ---
$ws.Cells.Item(169,1).Value = 11
$ws.Cells.Item(169,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(169,3).Value = "Bíobío"
$ws.Cells.Item(169,4).Value = 44250
$ws.Cells.Item(169,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(169,5).Value = 8
$ws.Cells.Item(169,6).Value = "Fruta"
$ws.Cells.Item(169,7).Value = 100101
$ws.Cells.Item(169,8).Value = "Berries"
$ws.Cells.Item(169,9).Value = 100112025
$ws.Cells.Item(169,10).Value = "Frutilla"
$ws.Cells.Item(169,11).Value = "Sin especificar"
$ws.Cells.Item(169,12).Value = "Primera"
$ws.Cells.Item(169,13).Value = 200
$ws.Cells.Item(169,14).Value = 8000
$ws.Cells.Item(169,15).Value = 8000
$ws.Cells.Item(169,16).Value = 8000
$ws.Cells.Item(169,17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(169,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(169,19).Value = 1143
$ws.Cells.Item(169,20).Value = 7

# --- New row 170 (appended at the end, mirrors old row 168's date) ---
$ws.Cells.Item(170,1).Value = 11
$ws.Cells.Item(170,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(170,3).Value = "Bíobío"
$ws.Cells.Item(170,4).Value = 44250
$ws.Cells.Item(170,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(170,5).Value = 8
$ws.Cells.Item(170,6).Value = "Fruta"
$ws.Cells.Item(170,7).Value = 100101
$ws.Cells.Item(170,8).Value = "Berries"
$ws.Cells.Item(170,9).Value = 100112025
$ws.Cells.Item(170,10).Value = "Frutilla"
$ws.Cells.Item(170,11).Value = "Sin especificar"
$ws.Cells.Item(170,12).Value = "Segunda"
$ws.Cells.Item(170,13).Value = 50
$ws.Cells.Item(170,14).Value = 6000
$ws.Cells.Item(170,15).Value = 6000
$ws.Cells.Item(170,16).Value = 6000
$ws.Cells.Item(170,17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(170,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(170,19).Value = 857
$ws.Cells.Item(170,20).Value = 7
